$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad" date) for rows 2-9 from 45212 to 45221
for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
